$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 434, pushing the existing rows
# 434-516 down to 435-517 (dimension grows from R516 to R517).
$ws.Rows.Item(434).Insert()

# Populate the newly inserted row 434 with the new record.
$ws.Range("A434").Value = 3
$ws.Range("B434").Value = "Femacal de La Calera"
$ws.Range("C434").Value = "Coquimbo"
$ws.Range("D434").Value = Get-Date -Year 2023 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("E434").Value = 5
$ws.Range("F434").Value = 100112043
$ws.Range("G434").Value = "Pepino ensalada"
$ws.Range("H434").Value = "Sin especificar"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 130
$ws.Range("K434").Value = 8500
$ws.Range("L434").Value = 9000
$ws.Range("M434").Value = 8769
$ws.Range("N434").Value = "$/caja 60 unidades"
$ws.Range("O434").Value = "Región de Arica y Parinacota"
$ws.Range("P434").Value = 146
$ws.Range("Q434").Value = 60
$ws.Range("R434").Value = "Hortaliza"
